$d = $word.ActiveDocument

$replacements = @(
    @{ old = "935×5=4675"; new = "488×7=3416" },
    @{ old = "157×8=1256"; new = "489×2=978" },
    @{ old = "769×7=5383"; new = "359×9=3231" },
    @{ old = "177×6=1062"; new = "496×6=2976" },
    @{ old = "611×2=1222"; new = "415×8=3320" },
    @{ old = "498×9=4482"; new = "616×9=5544" },
    @{ old = "458×6=2748"; new = "750×4=3000" },
    @{ old = "892×8=7136"; new = "391×9=3519" },
    @{ old = "554×8=4432"; new = "823×9=7407" },
    @{ old = "222×9=1998"; new = "787×8=6296" },
    @{ old = "852×5=4260"; new = "776×7=5432" },
    @{ old = "343×3=1029"; new = "294×2=588" },
    @{ old = "773×4=3092"; new = "338×3=1014" },
    @{ old = "699×3=2097"; new = "515×8=4120" },
    @{ old = "289×6=1734"; new = "736×9=6624" },
    @{ old = "603×8=4824"; new = "992×7=6944" },
    @{ old = "727×3=2181"; new = "661×6=3966" },
    @{ old = "952×3=2856"; new = "387×5=1935" },
    @{ old = "636×2=1272"; new = "313×3=939" },
    @{ old = "414×8=3312"; new = "694×5=3470" },
    @{ old = "627×6=3762"; new = "543×9=4887" },
    @{ old = "384×3=1152"; new = "390×6=2340" },
    @{ old = "892×2=1784"; new = "377×9=3393" },
    @{ old = "602×6=3612"; new = "634×5=3170" },
    @{ old = "288×2=576"; new = "987×3=2961" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
